# Apply updated probability values to team matrix (games pulled March 7)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2266187050359712
$ws.Range("C2").Value = 0.4496402877697842
$ws.Range("J2").Value = 0.01798561151079137
$ws.Range("P2").Value = 0.2014388489208633
$ws.Range("S2").Value = 0.1043165467625899
$ws.Range("B3").Value = 0.01587301587301587
$ws.Range("C3").Value = 0.03174603174603174
$ws.Range("J3").Value = 0.02380952380952381
$ws.Range("O3").Value = 0.007936507936507936
$ws.Range("P3").Value = 0.8015873015873016
$ws.Range("S3").Value = 0.119047619047619
$ws.Range("J4").Value = 0.0425531914893617
$ws.Range("O4").Value = 0.02127659574468085
$ws.Range("P4").Value = 0.6382978723404256
$ws.Range("S4").Value = 0.2978723404255319
$ws.Range("B6").Value = 0.07317073170731707
$ws.Range("D6").Value = 0.00975609756097561
$ws.Range("F6").Value = 0.05853658536585366
$ws.Range("J6").Value = 0.2585365853658537
$ws.Range("O6").Value = 0.01463414634146342
$ws.Range("Q6").Value = 0.1512195121951219
$ws.Range("R6").Value = 0.05853658536585366
$ws.Range("S6").Value = 0.375609756097561
$ws.Range("B7").Value = 0.1162790697674419
$ws.Range("D7").Value = 0.01550387596899225
$ws.Range("E7").Value = 0.007751937984496124
$ws.Range("F7").Value = 0.02713178294573643
$ws.Range("J7").Value = 0.1317829457364341
$ws.Range("O7").Value = 0.003875968992248062
$ws.Range("Q7").Value = 0.2170542635658915
$ws.Range("R7").Value = 0.05038759689922481
$ws.Range("S7").Value = 0.4302325581395349
$ws.Range("B8").Value = 0.06715063520871144
$ws.Range("D8").Value = 0.0235934664246824
$ws.Range("F8").Value = 0.05807622504537205
$ws.Range("J8").Value = 0.1161524500907441
$ws.Range("O8").Value = 0.01088929219600726
$ws.Range("Q8").Value = 0.1814882032667877
$ws.Range("R8").Value = 0.0852994555353902
$ws.Range("S8").Value = 0.4573502722323049
$ws.Range("B9").Value = 0.06622516556291391
$ws.Range("F9").Value = 0.07947019867549669
$ws.Range("J9").Value = 0.1125827814569536
$ws.Range("O9").Value = 0.01324503311258278
$ws.Range("Q9").Value = 0.1920529801324503
$ws.Range("R9").Value = 0.05960264900662252
$ws.Range("S9").Value = 0.4768211920529801
$ws.Range("B10").Value = 0.09090909090909091
$ws.Range("D10").Value = 0.02408702408702409
$ws.Range("F10").Value = 0.06604506604506605
$ws.Range("J10").Value = 0.1134421134421134
$ws.Range("O10").Value = 0.008547008547008548
$ws.Range("Q10").Value = 0.2362082362082362
$ws.Range("R10").Value = 0.06526806526806526
$ws.Range("S10").Value = 0.3954933954933955
$ws.Range("G11").Value = 0.1472222222222222
$ws.Range("J11").Value = 0.06666666666666667
$ws.Range("K11").Value = 0.2
$ws.Range("L11").Value = 0.5666666666666667
$ws.Range("S11").Value = 0.01944444444444444
$ws.Range("G12").Value = 0.7926267281105991
$ws.Range("J12").Value = 0.1152073732718894
$ws.Range("L12").Value = 0.04147465437788019
$ws.Range("S12").Value = 0.05069124423963134
$ws.Range("G13").Value = 0.6666666666666666
$ws.Range("J13").Value = 0.2666666666666667
$ws.Range("S13").Value = 0.06666666666666667
$ws.Range("F15").Value = 0.00975609756097561
$ws.Range("H15").Value = 0.2195121951219512
$ws.Range("I15").Value = 0.05365853658536585
$ws.Range("J15").Value = 0.3219512195121951
$ws.Range("K15").Value = 0.0975609756097561
$ws.Range("M15").Value = 0.01463414634146342
$ws.Range("O15").Value = 0.06341463414634146
$ws.Range("S15").Value = 0.2195121951219512
$ws.Range("F16").Value = 0.0213903743315508
$ws.Range("H16").Value = 0.1657754010695187
$ws.Range("I16").Value = 0.0481283422459893
$ws.Range("J16").Value = 0.3636363636363636
$ws.Range("K16").Value = 0.1764705882352941
$ws.Range("M16").Value = 0.0160427807486631
$ws.Range("O16").Value = 0.0213903743315508
$ws.Range("S16").Value = 0.1871657754010695
$ws.Range("F17").Value = 0.01740812379110251
$ws.Range("H17").Value = 0.2147001934235977
$ws.Range("I17").Value = 0.07156673114119923
$ws.Range("J17").Value = 0.4003868471953578
$ws.Range("K17").Value = 0.1083172147001934
$ws.Range("M17").Value = 0.01740812379110251
$ws.Range("O17").Value = 0.04642166344294004
$ws.Range("S17").Value = 0.1237911025145068
$ws.Range("F18").Value = 0.01818181818181818
$ws.Range("H18").Value = 0.2
$ws.Range("I18").Value = 0.103030303030303
$ws.Range("J18").Value = 0.4424242424242424
$ws.Range("K18").Value = 0.09090909090909091
$ws.Range("M18").Value = 0.0303030303030303
$ws.Range("N18").Value = 0.006060606060606061
$ws.Range("O18").Value = 0.04242424242424243
$ws.Range("S18").Value = 0.06666666666666667
$ws.Range("F19").Value = 0.01069900142653352
$ws.Range("H19").Value = 0.2417974322396576
$ws.Range("I19").Value = 0.05563480741797432
$ws.Range("J19").Value = 0.3573466476462197
$ws.Range("K19").Value = 0.1141226818830243
$ws.Range("M19").Value = 0.02995720399429386
$ws.Range("N19").Value = 0.0007132667617689016
$ws.Range("O19").Value = 0.07203994293865906
$ws.Range("S19").Value = 0.1176890156918688
